$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 45, shifting rows 45:174(old 173) down to 46:174
$ws.Rows.Item(45).Insert()

# Fill in the new row 45 with the new data point
$ws.Range("A45").Value = 11
$ws.Range("B45").Value = "Vega Monumental Concepción"
$ws.Range("C45").Value = "Bíobío"
$ws.Range("D45").Value = 44708
$ws.Range("D45").NumberFormat = $ws.Range("D46").NumberFormat
$ws.Range("E45").Value = 8
$ws.Range("F45").Value = 100112003
$ws.Range("G45").Value = "Ajo"
$ws.Range("H45").Value = "Chino"
$ws.Range("I45").Value = "Primera"
$ws.Range("J45").Value = 400
$ws.Range("K45").Value = 17000
$ws.Range("L45").Value = 18000
$ws.Range("M45").Value = 17500
$ws.Range("N45").Value = "$/caja 10 kilos"
$ws.Range("O45").Value = "China"
$ws.Range("P45").Value = 1750
$ws.Range("Q45").Value = 10
$ws.Range("R45").Value = "Hortaliza"
